$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the standalone "Meta description" paragraph that currently
#    follows the title (H1) paragraph.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------
# 2) Insert a new bold paragraph (re-using the page title text) right
#    before the final "Prompt: ..." paragraph.
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$pPrev = $d.Paragraphs.Item($n - 1)
$insertPoint = $d.Range($pPrev.Range.End, $pPrev.Range.End)
$insertPoint.InsertAfter("Play Desert Cats for Free - A Review of Ancient Egyptian Slot Game`r")

$pNew = $d.Paragraphs.Item($n)
$newTextRange = $d.Range($pNew.Range.Start, $pNew.Range.End - 1)
$newTextRange.Font.Bold = 1

# ---------------------------------------------------------------------
# 3) Replace the old "Prompt: ..." image-generation text (now the very
#    last paragraph) with the meta-description copy, keeping the
#    existing italic run formatting intact.
# ---------------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLast.Range.Find.Execute(
    "Prompt: Create a cartoon-style feature image for the game Desert Cats that features a happy Maya warrior with glasses. The image should be bright and colorful, with the happy Maya warrior front and center. The warrior should be depicted with a big grin on their face, wearing a colorful outfit and holding a scepter in one hand. They should also be wearing glasses, which should be prominently featured. In the background, there should be a desert landscape with palm trees and pyramids in the distance. The golden cats from the game should be incorporated into the image, lounging around the warrior's feet or perched on the scepter. The overall vibe of the image should be fun and playful, capturing the excitement of playing Desert Cats.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of the slot game Desert Cats with an Ancient Egyptian theme, featuring two Wild symbols and a Jackpot. Play for free here.",
    2)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
